# Add the new "Deleted Service doesn't affect the linked Service version"
# bug-report entry as row 14 of the BugReport tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title       = "Deleted Service doesn’t affect the linked Service version"
$description = "API under test: DELETE  /v1/services/{svcId}`nSteps to reproduce: `n1. Create a Service and the service version`n2. List the service version and see the service version fecthed successfully`n3. Delete the service`nSee that the deleted service has no impct on the linked service versions"
$apiArea     = "Delete Service`nQuestion/Clarification"

# ---------------------------------------------------------------------------
# Values
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = $title
$ws.Cells.Item(14, 3).Value = $description
$ws.Cells.Item(14, 4).Value = $apiArea

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------
# Column C reuses the sheet's standard "full border + wrap text" look, so copy
# it straight from the cell above rather than re-building it from scratch.
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# Columns A, B and D get a lighter border: thin left/right rule only, no
# top/bottom. Build it up evenly across all three cells (left edge on all of
# them, then right edge on all of them) so they end up sharing one style.
$cellA = $ws.Cells.Item(14, 1)
$cellB = $ws.Cells.Item(14, 2)
$cellD = $ws.Cells.Item(14, 4)

$cellA.Borders.Item(7).LineStyle = 1
$cellB.Borders.Item(7).LineStyle = 1
$cellD.Borders.Item(7).LineStyle = 1

$cellA.Borders.Item(10).LineStyle = 1
$cellB.Borders.Item(10).LineStyle = 1
$cellD.Borders.Item(10).LineStyle = 1

# D also wraps its text (title/question-clarification stack on two lines).
$cellD.WrapText = $true

# Row height to fit the wrapped description text.
$ws.Rows.Item(14).RowHeight = 99.75

# ---------------------------------------------------------------------------
# View state - land on the newly-added row like the author did after typing it
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
